$wb = $excel.ActiveWorkbook

# The same update needs to be applied to both the "展览" and "全部类型"
# worksheets, which contain duplicated data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F7").Value = 2363
    $ws.Range("F14").Value = 2
    $ws.Range("F15").Value = 973
}
